$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (CONTROLES -> Control); this also updates the
# _xlnm._FilterDatabase defined name reference automatically.
$ws.Name = "Control"

# Move the active selection on the sheet to D10 (was B1 with sqref A1:B1048576)
$ws.Range("D10").Select()
